# "gui AC ngay ko can qua 30p" - refresh the power/AC alarm report with the
# latest 06/05/2025 extract: rows 2-3 get new sites/faults, and five more
# rows (4-8) are appended for the rest of the batch (incl. SITE_OOS_BY_POWER
# alarms and the two freshly-widened "note" columns B/E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen columns B, C, E so the longer site/ticket names are readable ---
# (Excel's ColumnWidth setter only accepts whole-pixel steps, so we feed it
# the value whose resulting stored width lands closest to the target.)
$ws.Columns("B").ColumnWidth = 33.75
$ws.Columns("C").ColumnWidth = 18.75
$ws.Columns("E").ColumnWidth = 83.75

# --- Make room for rows 4-8 and inherit row 2's cell formatting (borders/
# alignment/number format) onto them before writing any values ---
$ws.Range("A2:H2").Copy() | Out-Null
$ws.Range("A4:H8").PasteSpecial(-4122) | Out-Null

# Row 2: SR_BVI027M_HNI
$ws.Range("A2").Value = "SR_BVI027M_HNI"
$ws.Range("B2").Value = "Thon-Lien-Tong-BVI_HNI"
$ws.Range("C2").Value = "POWER_AC_EAS"
$ws.Range("D2").Value = "06/05/2025 15:42:22"
$ws.Range("E2").ClearContents() | Out-Null
$ws.Range("F2").Value = "Trạm viễn thông loại 1"
$ws.Range("G2").Value = "Ba Vì"
$ws.Range("H2").Value = 1.08

# Row 3: UL_BVI133M_HNI
$ws.Range("A3").Value = "UL_BVI133M_HNI"
$ws.Range("B3").Value = "TONG-BAT-THON-TONG-LENH-BVI_HNI"
$ws.Range("C3").Value = "POWER_AC_EAS"
$ws.Range("D3").Value = "06/05/2025 15:37:25"
$ws.Range("E3").ClearContents() | Out-Null
$ws.Range("F3").Value = "Trạm viễn thông loại 3"
$ws.Range("G3").Value = "Ba Vì"
$ws.Range("H3").Value = 1.16

# Row 4: UL_DPG058M_HNI (new)
$ws.Range("A4").Value = "UL_DPG058M_HNI"
$ws.Range("B4").Value = "KCN-CAU-GAO-DPG_HNI"
$ws.Range("C4").Value = "POWER_AC_EAS"
$ws.Range("D4").Value = "06/05/2025 14:30:14"
$ws.Range("E4").ClearContents() | Out-Null
$ws.Range("F4").Value = "Trạm viễn thông loại 3"
$ws.Range("G4").Value = "Đan Phượng"
$ws.Range("H4").Value = 2.28

# Row 5: 2G_BVI010M_HNI (new)
$ws.Range("A5").Value = "2G_BVI010M_HNI"
$ws.Range("B5").Value = "Cam-Thuong-Thon-Van-Minh-BVI_HNI"
$ws.Range("C5").Value = "SITE_OOS_BY_POWER"
$ws.Range("D5").Value = "06/05/2025 13:42:12"
$ws.Range("E5").Value = "184205 - VTHN TĐML - HNI dựng lại côt thay cáp - 4 - thainh1 - 06/05/2025 13:56:39"
$ws.Range("F5").Value = "Trạm viễn thông loại 2"
$ws.Range("G5").Value = "Ba Vì"
$ws.Range("H5").Value = 3.08

# Row 6: 3G_BVI010M_HNI (new)
$ws.Range("A6").Value = "3G_BVI010M_HNI"
$ws.Range("B6").Value = "Cam-Thuong-Thon-Van-Minh-BVI_HNI"
$ws.Range("C6").Value = "SITE_OOS_BY_POWER"
$ws.Range("D6").Value = "06/05/2025 13:41:56"
$ws.Range("E6").Value = "184205 - VTHN TĐML - HNI dựng lại côt thay cáp - 4 - thainh1 - 06/05/2025 13:56:40"
$ws.Range("F6").Value = "Trạm viễn thông loại 2"
$ws.Range("G6").Value = "Ba Vì"
$ws.Range("H6").Value = 3.08

# Row 7: 4G-BVI010M-HNI (new)
$ws.Range("A7").Value = "4G-BVI010M-HNI"
$ws.Range("B7").Value = "Cam-Thuong-Thon-Van-Minh-BVI_HNI"
$ws.Range("C7").Value = "SITE_OOS_BY_POWER"
$ws.Range("D7").Value = "06/05/2025 13:41:51"
$ws.Range("E7").Value = "184205 - VTHN TĐML - HNI dựng lại côt thay cáp - 4 - thainh1 - 06/05/2025 13:56:40"
$ws.Range("F7").Value = "Trạm viễn thông loại 2"
$ws.Range("G7").Value = "Ba Vì"
$ws.Range("H7").Value = 3.08

# Row 8: SR_BVI010M_HNI (new)
$ws.Range("A8").Value = "SR_BVI010M_HNI"
$ws.Range("B8").Value = "Cam-Thuong-Thon-Van-Minh-BVI_HNI"
$ws.Range("C8").Value = "POWER_AC_EAS"
$ws.Range("D8").Value = "06/05/2025 06:18:28"
$ws.Range("E8").Value = "Mất nguồn AC - 1 - huongvl1 - 06/05/2025 10:23:32"
$ws.Range("F8").Value = "Trạm viễn thông loại 2"
$ws.Range("G8").Value = "Ba Vì"
$ws.Range("H8").Value = 10.47
